# Update the roster table on Sheet1.
# The underlying data was reshuffled: the six players that used to sit at
# the bottom of the list (De'Aaron Fox ... Kyle Kuzma) now appear right
# after the header row, the remaining players keep their relative order,
# and the final row (previously Goga Bitadze / Orlando Magic) is replaced
# with Nick Richards / Phoenix Suns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("De'Aaron Fox",  "PG",          "Sacramento Kings"),
    @("Shaedon Sharpe", "SG,SF",       "Portland Trail Blazers"),
    @("DeMar DeRozan",  "SF,PF",       "Sacramento Kings"),
    @("Kyle Kuzma",     "PF",          "Washington Wizards"),
    @("Evan Mobley",    "PF,C",        "Cleveland Cavaliers"),
    @("Ja Morant",      "PG",          "Memphis Grizzlies"),
    @("Gradey Dick",    "SG,SF",       "Toronto Raptors"),
    @("Brook Lopez",    "C",           "Milwaukee Bucks"),
    @("Mikal Bridges",  "SG,SF,PF",    "New York Knicks"),
    @("Nikola Vucevic", "PF,C",        "Chicago Bulls"),
    @("Klay Thompson",  "SG,SF",       "Dallas Mavericks"),
    @("Josh Giddey",    "PG,SG,SF",    "Chicago Bulls"),
    @("Tyler Herro",    "PG,SG",       "Miami Heat"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Miles Bridges",  "SF,PF",       "Charlotte Hornets"),
    @("Luka Doncic",    "PG,SG",       "Dallas Mavericks"),
    @("Nick Richards",  "C",           "Phoenix Suns")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
